# Slide 5 contains the planning table (shape "Inhaltsplatzhalter 5").
# The table is widened (second column gets more room) and three of the
# "Termin" date cells get a "Bis " ("until") prefix added to their text.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

# Widen the second grid column (Ramón Wilhelm's column) from 162pt to
# 223.77pt (2057400 EMU -> 2841879 EMU). Resizing a single column this
# way leaves the other three columns' widths untouched and grows the
# overall table/graphic-frame width accordingly (648pt -> 709.77pt,
# 8229600 EMU -> 9014079 EMU).
$tbl.Columns.Item(2).Width = 223.77

# Prefix the three date-only cells in the first column with "Bis ".
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Bis 12.07.2015"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Bis 30.08.2015"
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Bis 02.10.2015"

# The extra text plus the re-flowed column widen the table's rendered
# (rounded) height from 317.6pt to 296pt in absolute terms is not how
# PowerPoint actually grows it -- the table's bounding box grows to
# 296pt (4033520 EMU -> 3759200 EMU) because row 4 (the row with the
# most wrapped text) now needs more vertical room. Grow just that row;
# this keeps every other row's height exactly as it was.
$tbl.Rows.Item(4).Height = 179.2
